# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to Sheets per commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 612
$ws.Range("I9").Value = 549.25
$ws.Range("J9").Value = 674.75
$ws.Range("K9").Value = 549.25
$ws.Range("L9").Value = 674.75
$ws.Range("M9").Value = -380.25
$ws.Range("N9").Value = -1012.75
$ws.Range("H18").Value = 312.57144
$ws.Range("J18").Value = 260
$ws.Range("L18").Value = 260
$ws.Range("N18").Value = -828
$ws.Range("H40").Value = 17807.5
$ws.Range("I40").Value = 16211.5625
$ws.Range("J40").Value = 20999.375
$ws.Range("K40").Value = 16211.5625
$ws.Range("L40").Value = 20999.375
$ws.Range("M40").Value = -16036.5625
$ws.Range("N40").Value = -21349.375
$ws.Range("H53").Value = 896.1818
$ws.Range("J53").Value = 485.8
$ws.Range("L53").Value = 485.8
$ws.Range("N53").Value = -1759.8
$ws.Range("H96").Value = 12095.714
$ws.Range("I96").Value = 982.3
$ws.Range("J96").Value = 39879.25
$ws.Range("K96").Value = 2946.9
$ws.Range("L96").Value = 119637.75
$ws.Range("M96").Value = -1573.9
$ws.Range("N96").Value = -122383.75
$ws.Range("H98").Value = 2336.3333
$ws.Range("I98").Value = 2336.3333
$ws.Range("K98").Value = 2336.3333
$ws.Range("M98").Value = -838.3332999999998
$ws.Range("H107").Value = 611.625
$ws.Range("I107").Value = 276.07693
$ws.Range("J107").Value = 2065.6667
$ws.Range("K107").Value = 276.07693
$ws.Range("L107").Value = 2065.6667
$ws.Range("M107").Value = 1643.92307
$ws.Range("N107").Value = -5905.6667
$ws.Range("H113").Value = 3598.923
$ws.Range("I113").Value = 3255.1428
$ws.Range("K113").Value = 3255.1428
$ws.Range("M113").Value = -1.142800000000079
$ws.Range("H116").Value = 14714447
$ws.Range("I116").Value = 25008886
$ws.Range("J116").Value = 8103.9287
$ws.Range("K116").Value = 25008886
$ws.Range("L116").Value = 8103.9287
$ws.Range("M116").Value = -25005444
$ws.Range("N116").Value = -14987.9287
$ws.Range("H119").Value = 3266.6667
$ws.Range("J119").Value = 3266.6667
$ws.Range("L119").Value = 9800.000100000001
$ws.Range("N119").Value = -19476.0001
$ws.Range("H122").Value = 2336.3333
$ws.Range("I122").Value = 2336.3333
$ws.Range("K122").Value = 7008.999899999999
$ws.Range("M122").Value = -4558.999899999999
$ws.Range("H127").Value = 2815.1765
$ws.Range("I127").Value = 1128
$ws.Range("J127").Value = 3334.3076
$ws.Range("K127").Value = 3384
$ws.Range("L127").Value = 10002.9228
$ws.Range("M127").Value = 1576
$ws.Range("N127").Value = -19922.9228
$ws.Range("H132").Value = 174620
$ws.Range("I132").Value = 280818.34
$ws.Range("J132").Value = 22248.479
$ws.Range("K132").Value = 842455.02
$ws.Range("L132").Value = 66745.43700000001
$ws.Range("M132").Value = -839925.02
$ws.Range("N132").Value = -71805.43700000001
$ws.Range("H135").Value = 2531.9167
$ws.Range("I135").Value = 1164.4231
$ws.Range("K135").Value = 10479.8079
$ws.Range("M135").Value = -7944.8079
$ws.Range("H137").Value = 375094.88
$ws.Range("I137").Value = 560670.4
$ws.Range("K137").Value = 1682011.2
$ws.Range("M137").Value = -1679461.2
$ws.Range("H138").Value = 4233.3076
$ws.Range("I138").Value = 972.55
$ws.Range("J138").Value = 6790.7646
$ws.Range("K138").Value = 2917.65
$ws.Range("L138").Value = 20372.2938
$ws.Range("M138").Value = 2222.35
$ws.Range("N138").Value = -30652.2938
$ws.Range("H141").Value = 3663.3635
$ws.Range("I141").Value = 3481.6785
$ws.Range("J141").Value = 4680.8
$ws.Range("K141").Value = 10445.0355
$ws.Range("L141").Value = 14042.4
$ws.Range("M141").Value = -5265.0355
$ws.Range("N141").Value = -24402.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 56811.11
$ws.Range("I2").Value = 56811.11
$ws.Range("K2").Value = 56811.11
$ws.Range("M2").Value = -56698.11
$ws.Range("H32").Value = 12249.561
$ws.Range("I32").Value = 12346.612
$ws.Range("K32").Value = 12346.612
$ws.Range("M32").Value = -12059.612
$ws.Range("H61").Value = 4275.885
$ws.Range("I61").Value = 2704.5625
$ws.Range("K61").Value = 2704.5625
$ws.Range("M61").Value = -2492.5625
$ws.Range("H63").Value = 933
$ws.Range("I63").Value = 899.5
$ws.Range("J63").Value = 1000
$ws.Range("K63").Value = 899.5
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = -213.5
$ws.Range("N63").Value = -2372
$ws.Range("H66").Value = 933
$ws.Range("I66").Value = 899.5
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 4497.5
$ws.Range("L66").Value = 5000
$ws.Range("M66").Value = -1065.5
$ws.Range("N66").Value = -11864
$ws.Range("H74").Value = 1771.9166
$ws.Range("I74").Value = 1187.7142
$ws.Range("J74").Value = 2589.8
$ws.Range("K74").Value = 1187.7142
$ws.Range("L74").Value = 2589.8
$ws.Range("M74").Value = -313.7141999999999
$ws.Range("N74").Value = -4337.8
$ws.Range("H77").Value = 1771.9166
$ws.Range("I77").Value = 1187.7142
$ws.Range("J77").Value = 2589.8
$ws.Range("K77").Value = 5938.571
$ws.Range("L77").Value = 12949
$ws.Range("M77").Value = -1570.571
$ws.Range("N77").Value = -21685
$ws.Range("H88").Value = 63388
$ws.Range("I88").Value = 981
$ws.Range("J88").Value = 167399.67
$ws.Range("K88").Value = 981
$ws.Range("L88").Value = 167399.67
$ws.Range("M88").Value = -575
$ws.Range("N88").Value = -168211.67
$ws.Range("H91").Value = 63388
$ws.Range("I91").Value = 981
$ws.Range("J91").Value = 167399.67
$ws.Range("K91").Value = 981
$ws.Range("L91").Value = 167399.67
$ws.Range("M91").Value = 423
$ws.Range("N91").Value = -170207.67
$ws.Range("H112").Value = 71399.60000000001
$ws.Range("J112").Value = 71399.60000000001
$ws.Range("L112").Value = 71399.60000000001
$ws.Range("N112").Value = -74353.60000000001
$ws.Range("H116").Value = 56811.11
$ws.Range("I116").Value = 56811.11
$ws.Range("K116").Value = 56811.11
$ws.Range("M116").Value = -54517.11
$ws.Range("H132").Value = 15178.682
$ws.Range("J132").Value = 6259.9
$ws.Range("L132").Value = 18779.7
$ws.Range("N132").Value = -23839.7
$ws.Range("H136").Value = 4275.885
$ws.Range("I136").Value = 2704.5625
$ws.Range("K136").Value = 8113.6875
$ws.Range("M136").Value = -5563.6875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 56811.11
$ws.Range("I3").Value = 56811.11
$ws.Range("K3").Value = 56811.11
$ws.Range("M3").Value = -56697.11
$ws.Range("H80").Value = 893.2143
$ws.Range("I80").Value = 752.3333
$ws.Range("K80").Value = 752.3333
$ws.Range("M80").Value = 245.6667
$ws.Range("H83").Value = 893.2143
$ws.Range("I83").Value = 752.3333
$ws.Range("K83").Value = 3761.6665
$ws.Range("M83").Value = 1230.3335
$ws.Range("H105").Value = 3946.9167
$ws.Range("I105").Value = 3415.125
$ws.Range("K105").Value = 3415.125
$ws.Range("M105").Value = -1668.125
$ws.Range("H107").Value = 2322.0232
$ws.Range("I107").Value = 1843.0646
$ws.Range("K107").Value = 1843.0646
$ws.Range("M107").Value = 76.93540000000007
$ws.Range("H134").Value = 3440.3845
$ws.Range("I134").Value = 2550.5264
$ws.Range("K134").Value = 7651.5792
$ws.Range("M134").Value = -5116.5792

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4867.0356
$ws.Range("I31").Value = 3459.7646
$ws.Range("K31").Value = 3459.7646
$ws.Range("M31").Value = -3164.7646
$ws.Range("H34").Value = 4867.0356
$ws.Range("I34").Value = 3459.7646
$ws.Range("K34").Value = 3459.7646
$ws.Range("M34").Value = -3257.7646
$ws.Range("H58").Value = 2002402.4
$ws.Range("I58").Value = 2500
$ws.Range("J58").Value = 3335670.8
$ws.Range("K58").Value = 2500
$ws.Range("L58").Value = 3335670.8
$ws.Range("M58").Value = -2297
$ws.Range("N58").Value = -3336076.8
$ws.Range("H62").Value = 33575.625
$ws.Range("I62").Value = 27476.25
$ws.Range("J62").Value = 39675
$ws.Range("K62").Value = 27476.25
$ws.Range("L62").Value = 39675
$ws.Range("M62").Value = -26852.25
$ws.Range("N62").Value = -40923
$ws.Range("H65").Value = 33575.625
$ws.Range("I65").Value = 27476.25
$ws.Range("J65").Value = 39675
$ws.Range("K65").Value = 137381.25
$ws.Range("L65").Value = 198375
$ws.Range("M65").Value = -134261.25
$ws.Range("N65").Value = -204615
$ws.Range("H94").Value = 2695
$ws.Range("J94").Value = 3004.6
$ws.Range("L94").Value = 3004.6
$ws.Range("N94").Value = -3906.6
$ws.Range("H99").Value = 10639.6
$ws.Range("J99").Value = 10711.6875
$ws.Range("L99").Value = 10711.6875
$ws.Range("N99").Value = -13707.6875
$ws.Range("H100").Value = 85319.25
$ws.Range("J100").Value = 85319.25
$ws.Range("L100").Value = 85319.25
$ws.Range("N100").Value = -87483.25
$ws.Range("H107").Value = 808
$ws.Range("I107").Value = 646.25
$ws.Range("K107").Value = 646.25
$ws.Range("M107").Value = 1273.75
$ws.Range("H126").Value = 10639.6
$ws.Range("J126").Value = 10711.6875
$ws.Range("L126").Value = 32135.0625
$ws.Range("N126").Value = -37075.0625
$ws.Range("H132").Value = 6812837
$ws.Range("I132").Value = 7758434.5
$ws.Range("K132").Value = 23275303.5
$ws.Range("M132").Value = -23272773.5
$ws.Range("H133").Value = 82999.5
$ws.Range("J133").Value = 82999.5
$ws.Range("L133").Value = 82999.5
$ws.Range("N133").Value = -88059.5
$ws.Range("H134").Value = 2485.5806
$ws.Range("I134").Value = 2526.1738
$ws.Range("J134").Value = 2368.875
$ws.Range("K134").Value = 7578.5214
$ws.Range("L134").Value = 7106.625
$ws.Range("M134").Value = -5043.5214
$ws.Range("N134").Value = -12176.625
$ws.Range("H136").Value = 2002402.4
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 3335670.8
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 10007012.4
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -10012112.4
$ws.Range("H137").Value = 110000
$ws.Range("J137").Value = 110000
$ws.Range("L137").Value = 110000
$ws.Range("N137").Value = -120200

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 236.18182
$ws.Range("I2").Value = 186.66667
$ws.Range("J2").Value = 295.6
$ws.Range("K2").Value = 1120.00002
$ws.Range("L2").Value = 1773.6
$ws.Range("M2").Value = -1007.00002
$ws.Range("N2").Value = -1999.6
$ws.Range("H64").Value = 10707.3
$ws.Range("I64").Value = 1599.8
$ws.Range("J64").Value = 19814.8
$ws.Range("K64").Value = 4799.4
$ws.Range("L64").Value = 59444.39999999999
$ws.Range("M64").Value = -4529.4
$ws.Range("N64").Value = -59984.39999999999
$ws.Range("H67").Value = 10707.3
$ws.Range("I67").Value = 1599.8
$ws.Range("J67").Value = 19814.8
$ws.Range("K67").Value = 4799.4
$ws.Range("L67").Value = 59444.39999999999
$ws.Range("M67").Value = -3863.4
$ws.Range("N67").Value = -61316.39999999999
$ws.Range("H80").Value = 5998.227
$ws.Range("J80").Value = 5998.227
$ws.Range("L80").Value = 17994.681
$ws.Range("N80").Value = -19866.681
$ws.Range("H83").Value = 5998.227
$ws.Range("J83").Value = 5998.227
$ws.Range("L83").Value = 53984.043
$ws.Range("N83").Value = -63344.043
$ws.Range("H107").Value = 564.2143
$ws.Range("I107").Value = 297
$ws.Range("J107").Value = 712.6667
$ws.Range("K107").Value = 891
$ws.Range("L107").Value = 2138.0001
$ws.Range("M107").Value = 1029
$ws.Range("N107").Value = -5978.0001
$ws.Range("H121").Value = 1207.4286
$ws.Range("J121").Value = 1204
$ws.Range("L121").Value = 3612
$ws.Range("N121").Value = -6232
$ws.Range("H122").Value = 21254.445
$ws.Range("I122").Value = 95.5
$ws.Range("K122").Value = 859.5
$ws.Range("M122").Value = 1590.5
$ws.Range("H131").Value = 1957909.2
$ws.Range("I131").Value = 100935.5
$ws.Range("J131").Value = 2459794
$ws.Range("K131").Value = 302806.5
$ws.Range("L131").Value = 7379382
$ws.Range("M131").Value = -297766.5
$ws.Range("N131").Value = -7389462
$ws.Range("H132").Value = 5787.769
$ws.Range("J132").Value = 10583
$ws.Range("L132").Value = 95247
$ws.Range("N132").Value = -100307
$ws.Range("H137").Value = 141668260
$ws.Range("I137").Value = 150001900
$ws.Range("J137").Value = 100000000
$ws.Range("K137").Value = 450005700
$ws.Range("L137").Value = 300000000
$ws.Range("M137").Value = -450000600
$ws.Range("N137").Value = -300010200
$ws.Range("H140").Value = 8569.066000000001
$ws.Range("I140").Value = 2453.8
$ws.Range("K140").Value = 7361.400000000001
$ws.Range("M140").Value = -2181.400000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 102745.18
$ws.Range("I80").Value = 149286.86
$ws.Range("J80").Value = 21297.25
$ws.Range("K80").Value = 149286.86
$ws.Range("L80").Value = 21297.25
$ws.Range("M80").Value = -148288.86
$ws.Range("N80").Value = -23293.25
$ws.Range("H83").Value = 102745.18
$ws.Range("I83").Value = 149286.86
$ws.Range("J83").Value = 21297.25
$ws.Range("K83").Value = 746434.2999999999
$ws.Range("L83").Value = 106486.25
$ws.Range("M83").Value = -741442.2999999999
$ws.Range("N83").Value = -116470.25
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents() | Out-Null
$ws.Range("H132").Value = 604896.8
$ws.Range("I132").Value = 122584.586
$ws.Range("K132").Value = 367753.758
$ws.Range("M132").Value = -365223.758

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1054.1052
$ws.Range("I22").Value = 913.2308
$ws.Range("K22").Value = 913.2308
$ws.Range("M22").Value = -618.2308
$ws.Range("H27").Value = 1054.1052
$ws.Range("I27").Value = 913.2308
$ws.Range("K27").Value = 913.2308
$ws.Range("M27").Value = -806.2308
$ws.Range("H46").Value = 2817.2273
$ws.Range("I46").Value = 721.5
$ws.Range("K46").Value = 721.5
$ws.Range("M46").Value = -533.5
$ws.Range("H68").Value = 2941.96
$ws.Range("I68").Value = 2386
$ws.Range("J68").Value = 4702.5
$ws.Range("K68").Value = 2386
$ws.Range("L68").Value = 4702.5
$ws.Range("M68").Value = -1637
$ws.Range("N68").Value = -6200.5
$ws.Range("H71").Value = 2941.96
$ws.Range("I71").Value = 2386
$ws.Range("J71").Value = 4702.5
$ws.Range("K71").Value = 11930
$ws.Range("L71").Value = 23512.5
$ws.Range("M71").Value = -8186
$ws.Range("N71").Value = -31000.5
$ws.Range("H100").Value = 1373.6666
$ws.Range("I100").Value = 1373.6666
$ws.Range("K100").Value = 1373.6666
$ws.Range("M100").Value = -832.6666
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents() | Out-Null
$ws.Range("H132").Value = 2181.75
$ws.Range("I132").Value = 2210.9238
$ws.Range("J132").Value = 1846.25
$ws.Range("K132").Value = 6632.7714
$ws.Range("L132").Value = 5538.75
$ws.Range("M132").Value = -4102.7714
$ws.Range("N132").Value = -10598.75
$ws.Range("H136").Value = 2400.7285
$ws.Range("I136").Value = 1623.2963
$ws.Range("K136").Value = 4869.8889
$ws.Range("M136").Value = -2319.8889

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 112325.3
$ws.Range("J41").Value = 112325.3
$ws.Range("L41").Value = 112325.3
$ws.Range("N41").Value = -113105.3
$ws.Range("H87").Value = 60000
$ws.Range("I87").Value = 60000
$ws.Range("K87").Value = 60000
$ws.Range("M87").Value = -58752
$ws.Range("H90").Value = 60000
$ws.Range("I90").Value = 60000
$ws.Range("K90").Value = 180000
$ws.Range("M90").Value = -173760
$ws.Range("H100").Value = 3366.2415
$ws.Range("I100").Value = 3385.2856
$ws.Range("K100").Value = 6770.5712
$ws.Range("M100").Value = -6229.5712
$ws.Range("H113").Value = 1715.909
$ws.Range("I113").Value = 1003.25
$ws.Range("J113").Value = 3616.3333
$ws.Range("K113").Value = 3009.75
$ws.Range("L113").Value = 10848.9999
$ws.Range("M113").Value = -839.75
$ws.Range("N113").Value = -15188.9999
$ws.Range("H121").Value = 35000
$ws.Range("J121").Value = 35000
$ws.Range("L121").Value = 35000
$ws.Range("N121").Value = -38494
$ws.Range("H132").Value = 3571.9119
$ws.Range("I132").Value = 1121.5
$ws.Range("J132").Value = 11535.75
$ws.Range("K132").Value = 3364.5
$ws.Range("L132").Value = 34607.25
$ws.Range("M132").Value = -834.5
$ws.Range("N132").Value = -39667.25
$ws.Range("H136").Value = 9020.84
$ws.Range("I136").Value = 1924.3684
$ws.Range("K136").Value = 5773.1052
$ws.Range("M136").Value = -3223.1052
